$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Update the "gueltig ab:" (valid from) date value in cell B5
$ws.Range("B5").Value = "01.01.2024"
